$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently holds 7 data rows (rows 2-8, below the header in row 1).
# Two brand-new phone entries were added to the top of the list, pushing all
# the existing entries down by two rows (old row 2 -> new row 4, etc.).
$ws.Range("A2:C3").EntireRow.Insert()

# Row 4 (the former row 2) still carries the correct data-row style. Copy its
# formatting onto the two freshly inserted, currently blank rows so they look
# like the rest of the table instead of inheriting the header's style.
$ws.Range("A4:C4").Copy()
$ws.Range("A2:C3").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Force these cells to store their content as text (matching the rest of the
# column, which holds phone numbers, DDD codes and dates as plain strings)
# instead of letting Excel auto-convert them to numbers/dates.
$ws.Range("A2:C3").NumberFormat = "@"

# New row 2: +555491557534 / 54 / 2024-10-04
$ws.Range("A2").Value = "+555491557534"
$ws.Range("B2").Value = "54"
$ws.Range("C2").Value = "2024-10-04"

# New row 3: +556181971614 / 61 / 2024-10-01
$ws.Range("A3").Value = "+556181971614"
$ws.Range("B3").Value = "61"
$ws.Range("C3").Value = "2024-10-01"
